$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values on row 2
$ws.Range("D2").Value = 4
$ws.Range("F2").Value = -3
$ws.Range("H2").Value = 46

# Update selected cell
$ws.Range("C2").Select()
